# Junction_Flooding_463 edit: update simulation data (rows 2-5), remove row 6,
# and widen several data columns from 7 to 8 (raw OOXML width units).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Replace the simulation data in rows 2-5 with the new dataset ---
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH")

$row2 = @(45144.50694444445,9.173,6.534,2.853,20.581,14.857,6.617,20.037,11.575,4.738,6.004,8.497,8.644,2.821,7.519,9.989000000000001,7.075,2.196,0.759,107.045,20.658,6.941,12.856,7.392,1.396,11.913,6.131,5.786,6.628,8.785,2.215,18.053,3.621,8.676)
for ($i = 0; $i -lt $cols.Length; $i++) {
  $ws.Range($cols[$i] + "2").Value = $row2[$i]
}

$row3 = @(45144.51388888889,21.356,15.897,1.818,46.992,37.735,16.383,61.483,26.123,11.608,16.579,19.059,19.894,5.726,16.918,23.828,14.566,1.187,0.821,250.193,47.253,15.616,31.419,16.724,2.324,31.659,13.794,12.348,14.485,19.828,1.023,56.323,8.683,19.524)
for ($i = 0; $i -lt $cols.Length; $i++) {
  $ws.Range($cols[$i] + "3").Value = $row3[$i]
}

$row4 = @(45144.52083333334,0.292,0.132,0.718,0.834,0,0,10.154,0.531,0.27,0.354,0.277,0,0,0.376,0.6,0.572,0.6879999999999999,0.08,0,1.785,0.347,1.13,0.542,0.097,5.006,0.307,0.477,0.499,0.433,0.611,10.259,0.031,0.465)
for ($i = 0; $i -lt $cols.Length; $i++) {
  $ws.Range($cols[$i] + "4").Value = $row4[$i]
}

$row5 = @(45144.52777777778,11.38,8.49,0.9,24.94,20.07,9.220000000000001,31.58,13.91,6.06,9.33,9.949999999999999,10.34,2.79,9.02,12.65,7.69,0.62,0.37,130.41,25.03,8.33,16.62,8.859999999999999,1.21,15.97,7.36,6.59,7.74,10.58,0.46,28.4,4.62,10.41)
for ($i = 0; $i -lt $cols.Length; $i++) {
  $ws.Range($cols[$i] + "5").Value = $row5[$i]
}

# --- 2) Delete the now-unused row 6 (shrinks used range to A1:AH5) ---
$ws.Rows.Item(6).Delete()

# --- 3) Widen columns C, G, J, O, Q, V, X, AA, AB, AC from raw width 7 to 8 ---
# ColumnWidth (Excel character units) = raw OOXML width - 0.83 for this font/zoom,
# so 7.17 yields a saved <col width="8"/>. Columns.Item needs a numeric index here
# (passing the letter as a string throws a type-mismatch in this COM runtime).
$wideColsIdx = @(3,7,10,15,17,22,24,27,28,29)
foreach ($c in $wideColsIdx) {
  $ws.Columns.Item($c).ColumnWidth = 7.17
}

Write-Output "done"